$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (prevents Excel
# from auto-converting numeric-looking strings like "27.98" into numbers),
# then restore the default "Normal" style so no stray formatting is introduced.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "68.313.71"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
Set-TextValue "D3" "2.644.17"
$ws.Range("E3").Value = "  +0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
Set-TextValue "D5" "598.43"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
Set-TextValue "D6" "154.78"

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("E8").Value = "  -0.65%  "

# Row 9
Set-TextValue "D9" "2.643.23"
$ws.Range("E9").Value = "  +0.57%  "

# Row 10
Set-TextValue "D10" "0.146"
$ws.Range("E10").Value = "  +8.34%  "

# Row 11
$ws.Range("E11").Value = "  -0.59%  "

# Row 12
$ws.Range("E12").Value = "  +1.04%  "

# Row 13
$ws.Range("E13").Value = "  +2.15%  "

# Row 14
Set-TextValue "D14" "0.0000194"
$ws.Range("E14").Value = "  +3.00%  "

# Row 15
Set-TextValue "D15" "27.98"
$ws.Range("E15").Value = "  +1.41%  "

# Row 16
Set-TextValue "D16" "3.123.47"
$ws.Range("E16").Value = "  +0.56%  "

# Row 17
Set-TextValue "D17" "68.218.63"
$ws.Range("E17").Value = "  +0.59%  "

# Row 18
Set-TextValue "D18" "2.651.33"
$ws.Range("E18").Value = "  +0.92%  "

# Row 19
Set-TextValue "D19" "11.38"
$ws.Range("E19").Value = "  -0.53%  "

# Row 20
Set-TextValue "D20" "364.79"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21
Set-TextValue "D21" "7.47"
$ws.Range("E21").Value = "  +0.36%  "

# Row 22
$ws.Range("E22").Value = "  +3.30%  "

# Row 23
Set-TextValue "D23" "4.88"
$ws.Range("E23").Value = "  +1.75%  "

# Row 24
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D24" "2.06"
$ws.Range("E24").Value = "  -0.84%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D25" "75.53"
$ws.Range("E25").Value = "  +4.74%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
Set-TextValue "D27" "9.78"
$ws.Range("E27").Value = "  -1.02%  "

# Row 28
$ws.Range("E28").Value = "  +2.32%  "

# Row 30
$ws.Range("E30").Value = "  -0.03%  "

# Row 31
Set-TextValue "D31" "563.48"
$ws.Range("E31").Value = "  -2.24%  "

# Row 32
Set-TextValue "D32" "8.09"
$ws.Range("E32").Value = "  +2.10%  "

# Row 33
$ws.Range("E33").Value = "  +0.61%  "

# Row 34
$ws.Range("E34").Value = "  +1.06%  "

# Row 35
Set-TextValue "D35" "0.131"
$ws.Range("E35").Value = "  +3.12%  "

# Row 36
Set-TextValue "D36" "0.999"

# Row 37
Set-TextValue "D37" "1.58"
$ws.Range("E37").Value = "  +3.48%  "

# Row 38
Set-TextValue "D38" "161.24"
$ws.Range("E38").Value = "  +1.75%  "

# Row 39
$ws.Range("E39").Value = "  +0.86%  "

# Row 40
Set-TextValue "D40" "0.376"
$ws.Range("E40").Value = "  +2.11%  "

# Row 41
$ws.Range("E41").Value = "  -0.29%  "

# Row 42
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
Set-TextValue "D43" "0.0₆0341"
$ws.Range("E43").Value = "  +2.87%  "

# Row 44
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("E45").Value = "  +2.10%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D46" "40.62"
$ws.Range("E46").Value = "  +1.17%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "1.00"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
Set-TextValue "D48" "156.36"
$ws.Range("E48").Value = "  +0.55%  "

# Row 49
$ws.Range("E49").Value = "  +1.97%  "

# Row 50
$ws.Range("E50").Value = "  +0.50%  "

# Row 51
Set-TextValue "D51" "21.86"
$ws.Range("E51").Value = "  -0.51%  "
